# Generate Report for Handback
# Updates the localization-status report to reflect that handback has
# completed for both the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdDisplay  = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5beccf39e6056bd4c6c851dafad873ed9a4f98e/e2e/a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.md"
$ffDisplay  = "ffffe662beed-9921-4904-9775-ac1f3056ca94.md"
$ffUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d5beccf39e6056bd4c6c851dafad873ed9a4f98e/e2e/ffffe662beed-9921-4904-9775-ac1f3056ca94.md"

# ---------------------------------------------------------------------
# 1. Update status text on the Overview sheet (zh-cn / de-de columns)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# 2. zh-cn sheet: Status, Latest Target File, Latest Handback File,
#    Latest Handback DateTime
# ---------------------------------------------------------------------
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = $mdDisplay
$wsZh.Range("I3").Value = $mdDisplay

$wsZh.Range("J2").Value = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.d2780beeddf58ea9ef886415d987fd9b265baee5.zh-cn.xlf"
$wsZh.Range("J3").Value = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.d2780beeddf58ea9ef886415d987fd9b265baee5.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-15 09:16:31"
$wsZh.Range("K3").Value = "2016-08-15 09:16:31"

# Rebuild hyperlinks in the same order the original report builder uses:
# A2, I2, A3, I3 -- this reproduces the relationship id numbering seen in
# the generated workbook (rId2..rId5).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ffUrl, "", "", $ffDisplay)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl, "", "", $mdDisplay)

# ---------------------------------------------------------------------
# 3. de-de sheet: Status, Latest Target File, Latest Handback File,
#    Latest Handback DateTime
# ---------------------------------------------------------------------
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = $mdDisplay
$wsDe.Range("I3").Value = $mdDisplay

$wsDe.Range("J2").Value = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.d2780beeddf58ea9ef886415d987fd9b265baee5.de-de.xlf"
$wsDe.Range("J3").Value = "a0f10a7c-6e84-4123-8ad9-5bd465a4c7ad.d2780beeddf58ea9ef886415d987fd9b265baee5.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-15 09:16:38"
$wsDe.Range("K3").Value = "2016-08-15 09:16:38"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ffUrl, "", "", $ffDisplay)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl, "", "", $mdDisplay)

# ---------------------------------------------------------------------
# 4. Column width adjustments (the Status / Latest Target / Latest
#    Handback File columns widen to accommodate the new text).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.98
$wsOverview.Columns.Item(6).ColumnWidth = 29.98

$wsZh.Columns.Item(3).ColumnWidth = 29.98
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsDe.Columns.Item(3).ColumnWidth = 29.98
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
